$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 29 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4371.6
$ws.Range("I40").Value = 3252.8333
$ws.Range("J40").Value = 6049.75
$ws.Range("K40").Value = 3252.8333
$ws.Range("L40").Value = 6049.75
$ws.Range("M40").Value = -3077.8333
$ws.Range("N40").Value = -6399.75
$ws.Range("H135").Value = 554.13043
$ws.Range("I135").Value = 474.45
$ws.Range("K135").Value = 4270.05
$ws.Range("M135").Value = -1735.05
$ws.Range("H137").Value = 1697.3
$ws.Range("I137").Value = 1282
$ws.Range("J137").Value = 2666.3333
$ws.Range("K137").Value = 3846
$ws.Range("L137").Value = 7998.999899999999
$ws.Range("M137").Value = -1296
$ws.Range("N137").Value = -13098.9999
$ws.Range("H138").Value = 1801.295
$ws.Range("J138").Value = 2333.658
$ws.Range("L138").Value = 7000.974
$ws.Range("N138").Value = -17280.974
$ws.Range("H141").Value = 3627.7778
$ws.Range("I141").Value = 4763.4
$ws.Range("J141").Value = 2208.25
$ws.Range("K141").Value = 14290.2
$ws.Range("L141").Value = 6624.75
$ws.Range("M141").Value = -9110.199999999999
$ws.Range("N141").Value = -16984.75

# --- Sheet ARM: 45 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7499.8335
$ws.Range("I45").Value = 4999.6665
$ws.Range("K45").Value = 4999.6665
$ws.Range("M45").Value = -4622.6665
$ws.Range("H61").Value = 8677.25
$ws.Range("I61").Value = 7774
$ws.Range("K61").Value = 7774
$ws.Range("M61").Value = -7562
$ws.Range("H63").Value = 5143.6665
$ws.Range("I63").Value = 3800
$ws.Range("J63").Value = 6487.3335
$ws.Range("K63").Value = 3800
$ws.Range("L63").Value = 6487.3335
$ws.Range("M63").Value = -3114
$ws.Range("N63").Value = -7859.3335
$ws.Range("H66").Value = 5143.6665
$ws.Range("I66").Value = 3800
$ws.Range("J66").Value = 6487.3335
$ws.Range("K66").Value = 19000
$ws.Range("L66").Value = 32436.6675
$ws.Range("M66").Value = -15568
$ws.Range("N66").Value = -39300.6675
$ws.Range("H74").Value = 2113.2698
$ws.Range("I74").Value = 1579.625
$ws.Range("K74").Value = 1579.625
$ws.Range("M74").Value = -705.625
$ws.Range("H77").Value = 2113.2698
$ws.Range("I77").Value = 1579.625
$ws.Range("K77").Value = 7898.125
$ws.Range("M77").Value = -3530.125
$ws.Range("H122").Value = 3098.261
$ws.Range("I122").Value = 2949.8667
$ws.Range("J122").Value = 3376.5
$ws.Range("K122").Value = 8849.6001
$ws.Range("L122").Value = 10129.5
$ws.Range("M122").Value = -6399.6001
$ws.Range("N122").Value = -15029.5
$ws.Range("H132").Value = 3565.5557
$ws.Range("I132").Value = 3010.1667
$ws.Range("K132").Value = 9030.500100000001
$ws.Range("M132").Value = -6500.500100000001
$ws.Range("H136").Value = 8677.25
$ws.Range("I136").Value = 7774
$ws.Range("K136").Value = 23322
$ws.Range("M136").Value = -20772

# --- Sheet BSM: 4 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2860.6667
$ws.Range("I54").Value = 2860.6667
$ws.Range("K54").Value = 2860.6667
$ws.Range("M54").Value = -2376.6667

# --- Sheet CRP: 36 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2733.25
$ws.Range("I35").Value = 2733.25
$ws.Range("K35").Value = 2733.25
$ws.Range("M35").Value = -2439.25
$ws.Range("H58").Value = 6698.875
$ws.Range("I58").Value = 3798.25
$ws.Range("K58").Value = 3798.25
$ws.Range("M58").Value = -3595.25
$ws.Range("H82").Value = 55000
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55722
$ws.Range("H85").Value = 55000
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57496
$ws.Range("H99").Value = 4245.3
$ws.Range("I99").Value = 4150.5
$ws.Range("K99").Value = 4150.5
$ws.Range("M99").Value = -2652.5
$ws.Range("H126").Value = 4245.3
$ws.Range("I126").Value = 4150.5
$ws.Range("K126").Value = 12451.5
$ws.Range("M126").Value = -9981.5
$ws.Range("H134").Value = 4842.2188
$ws.Range("I134").Value = 4317.48
$ws.Range("K134").Value = 12952.44
$ws.Range("M134").Value = -10417.44
$ws.Range("H136").Value = 6698.875
$ws.Range("I136").Value = 3798.25
$ws.Range("K136").Value = 11394.75
$ws.Range("M136").Value = -8844.75
$ws.Range("H141").Value = 32785.582
$ws.Range("J141").Value = 34997.832
$ws.Range("L141").Value = 34997.832
$ws.Range("N141").Value = -45357.832

# --- Sheet CUL: 35 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10025
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 150
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -38
$ws.Range("N3").Value = -60224
$ws.Range("H11").Value = 33333934
$ws.Range("I11").Value = 100000000
$ws.Range("K11").Value = 300000000
$ws.Range("M11").Value = -299999860
$ws.Range("H34").Value = 293.1111
$ws.Range("H39").Value = 1306.2858
$ws.Range("J39").Value = 1468.8
$ws.Range("L39").Value = 4406.4
$ws.Range("N39").Value = -4994.4
$ws.Range("H55").Value = 718.8
$ws.Range("I55").Value = 297.5
$ws.Range("J55").Value = 999.6667
$ws.Range("K55").Value = 892.5
$ws.Range("L55").Value = 2999.0001
$ws.Range("M55").Value = -715.5
$ws.Range("N55").Value = -3353.0001
$ws.Range("H81").Value = 27747.25
$ws.Range("J81").Value = 36663
$ws.Range("L81").Value = 109989
$ws.Range("N81").Value = -112235
$ws.Range("H84").Value = 27747.25
$ws.Range("J84").Value = 36663
$ws.Range("L84").Value = 329967
$ws.Range("N84").Value = -341199
$ws.Range("H129").Value = 18126946
$ws.Range("J129").Value = 1006672.6
$ws.Range("L129").Value = 3020017.8
$ws.Range("N129").Value = -3030017.8

# --- Sheet GSM: 43 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18400
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 18400
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18400
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -18940
$ws.Range("H73").Value = 18400
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 18400
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18400
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -20272
$ws.Range("H80").Value = 3719.3
$ws.Range("I80").Value = 4109.357
$ws.Range("J80").Value = 3378
$ws.Range("K80").Value = 4109.357
$ws.Range("L80").Value = 3378
$ws.Range("M80").Value = -3111.357
$ws.Range("N80").Value = -5374
$ws.Range("H83").Value = 3719.3
$ws.Range("I83").Value = 4109.357
$ws.Range("J83").Value = 3378
$ws.Range("K83").Value = 20546.785
$ws.Range("L83").Value = 16890
$ws.Range("M83").Value = -15554.785
$ws.Range("N83").Value = -26874
$ws.Range("H113").Value = 1583
$ws.Range("I113").Value = 1499.5
$ws.Range("K113").Value = 1499.5
$ws.Range("M113").Value = 670.5
$ws.Range("H122").Value = 4171.0713
$ws.Range("I122").Value = 3737.1667
$ws.Range("J122").Value = 6774.5
$ws.Range("K122").Value = 11211.5001
$ws.Range("L122").Value = 20323.5
$ws.Range("M122").Value = -8761.500100000001
$ws.Range("N122").Value = -25223.5
$ws.Range("H126").Value = 4211
$ws.Range("I126").Value = 3849.1875
$ws.Range("K126").Value = 11547.5625
$ws.Range("M126").Value = -9077.5625

# --- Sheet LTW: 18 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4499.6665
$ws.Range("I22").Value = 6000
$ws.Range("J22").Value = 3749.5
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 3749.5
$ws.Range("M22").Value = -5705
$ws.Range("N22").Value = -4339.5
$ws.Range("H27").Value = 4499.6665
$ws.Range("I27").Value = 6000
$ws.Range("J27").Value = 3749.5
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 3749.5
$ws.Range("M27").Value = -5893
$ws.Range("N27").Value = -3963.5
$ws.Range("H40").Value = 953.3333
$ws.Range("I40").Value = 953.3333
$ws.Range("K40").Value = 953.3333
$ws.Range("M40").Value = -817.3333

# --- Sheet WVR: 22 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5582.846
$ws.Range("I62").Value = 5019.7
$ws.Range("K62").Value = 5019.7
$ws.Range("M62").Value = -4395.7
$ws.Range("H65").Value = 5582.846
$ws.Range("I65").Value = 5019.7
$ws.Range("K65").Value = 25098.5
$ws.Range("M65").Value = -21978.5
$ws.Range("H100").Value = 1711.1538
$ws.Range("I100").Value = 1675.2222
$ws.Range("J100").Value = 1792
$ws.Range("K100").Value = 3350.4444
$ws.Range("L100").Value = 3584
$ws.Range("M100").Value = -2809.4444
$ws.Range("N100").Value = -4666
$ws.Range("H107").Value = 759.875
$ws.Range("I107").Value = 613.1667
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1839.5001
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = 80.49990000000003
$ws.Range("N107").Value = -7440
